# "Generate Report for Handback" -- localization-status.xlsx refresh
#
# The CI job re-ran and produced a fresh handback report:
#   - Status flips from "Ready for handoff" to "Handed back: in sync with en-US"
#     (this shared string is shown on the Overview sheet's zh-cn/de-de columns
#     as well as on each language sheet's Status column, since they all share
#     the same text).
#   - The "Latest Handback DateTime" for each language gets a newer timestamp.
#   - The stale "version mismatch" Error Detail message is cleared now that the
#     handback is in sync.
#   - A couple of columns are widened in the refreshed report layout.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# Update every cell that carries this value so the old string becomes fully
# unreferenced (Overview mirrors the per-language Status column).
$ws1.Range("E2").Value = "Handed back: in sync with en-US"
$ws1.Range("F2").Value = "Handed back: in sync with en-US"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"

# --- Latest Handback DateTime refresh ---
$ws2.Range("K2").Value = "2016-09-04 22:53:43"
$ws3.Range("K2").Value = "2016-09-04 22:53:51"

# --- Error Detail cleared (handback is now in sync, no version-mismatch error) ---
$ws2.Range("P2").Value = "'"
$ws3.Range("P2").Value = "'"

# --- Column width adjustments for the refreshed report ---
$ws1.Columns.Item(5).ColumnWidth = 29.166666666666668
$ws1.Columns.Item(6).ColumnWidth = 29.166666666666668

$ws2.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws2.Columns.Item(16).ColumnWidth = 12.833333333333334

$ws3.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws3.Columns.Item(16).ColumnWidth = 12.833333333333334
